$d = $word.ActiveDocument

# Position a collapsed range at the very end of the document body
# (i.e. right after the last paragraph, before the sectPr).
$endPos = $d.Content.End
$rng = $d.Range($endPos, $endPos)

# Build the raw WordprocessingML for the two paragraphs to append:
#   1) an empty paragraph (keeps the fr-FR paragraph-mark language)
#   2) a paragraph containing the text "cvbjk" (also fr-FR)
# Using Range.InsertXML lets us insert exactly this markup without the
# engine auto-materializing an empty run in the blank paragraph (which
# happens if we build it via Selection.TypeParagraph/TypeText instead).
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
      '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
          '<w:body>' + `
            '<w:p>' + `
              '<w:pPr><w:rPr><w:lang w:val="fr-FR"/></w:rPr></w:pPr>' + `
            '</w:p>' + `
            '<w:p>' + `
              '<w:pPr><w:rPr><w:lang w:val="fr-FR"/></w:rPr></w:pPr>' + `
              '<w:r><w:rPr><w:lang w:val="fr-FR"/></w:rPr><w:t>cvbjk</w:t></w:r>' + `
            '</w:p>' + `
          '</w:body>' + `
        '</w:document>' + `
      '</pkg:xmlData>' + `
    '</pkg:part>' + `
  '</pkg:package>'

$rng.InsertXML($xml)
